$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 251; everything from old row 251 down
# (through old row 315) shifts down to 252..316, matching the target
# dimension A1:R316.
$ws.Rows("251:251").Insert()

# Populate the newly inserted row 251 with the new weekly data point.
$ws.Range("A251").Value = 6
$ws.Range("B251").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C251").Value = "Metropolitana"
$ws.Range("D251").Value = 44511
$ws.Range("E251").Value = 13
$ws.Range("F251").Value = 100112039
$ws.Range("G251").Value = "Ciboulette"
$ws.Range("H251").Value = "Sin especificar"
$ws.Range("I251").Value = "Primera"
$ws.Range("J251").Value = 840
$ws.Range("K251").Value = 700
$ws.Range("L251").Value = 800
$ws.Range("M251").Value = 754
$ws.Range("N251").Value = "`$/docena de atados"
$ws.Range("O251").Value = "Región Metropolitana"
$ws.Range("P251").Value = 251
$ws.Range("Q251").Value = 3
$ws.Range("R251").Value = "Hortaliza"
